$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "date",
    "nea-ukcharity.bsky.social",
    "caneurope.bsky.social",
    "wateraid.bsky.social",
    "migrantsrights.bsky.social",
    "friends-earth.bsky.social",
    "samcardwell44.bsky.social",
    "greenpeace.eu",
    "wwfeu.bsky.social",
    "powertochange.org.uk",
    "thegreenregister.bsky.social",
    "endfuelpoverty.bsky.social",
    "commenergyengland.bsky.social",
    "extinctionrebellion.uk",
    "wwtworldwide.bsky.social",
    "bristolgreenparty.bsky.social",
    "warmthiswinter.bsky.social",
    "jrct.bsky.social",
    "ssencommunity.bsky.social",
    "localtrust.bsky.social",
    "wiltscouncil.bsky.social",
    "nationalgrid.bsky.social",
    "ofgem.bsky.social",
    "barnsleycouncil.bsky.social",
    "northsomersetc.bsky.social",
    "citizensadvice.bsky.social"
)

# Clear out the old row 2 ("date" sample value) - new layout is header-only.
$ws.Range("A2").Value = $null

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
}

$ws.Range("A2:Z2").Select()
